# Implemented with ArrayList<Integer>. Cut time by ~60%.
#
# This edit:
#  1. Renames the existing "Sheet1" to "v1.0".
#  2. Duplicates it as a new sheet "v1.1 ArrayList<Integer>" placed right
#     after "v1.0", which becomes the active sheet/tab.
#  3. Updates the new sheet's raw timing inputs (column A) and the
#     FORECAST.LINEAR ranges (column E) to reflect the smaller (6-row)
#     sample set used by the ArrayList<Integer> implementation, clearing
#     out the now-unused rows 8-9 of columns A/B.
#  4. Restores the selections / active-cell bookmarks left behind on each
#     tab, matching what was recorded when the workbook was last saved.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the original sheet -----------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "v1.0"

# --- 2. Duplicate it right after itself to create the v1.1 sheet ----------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "v1.1 ArrayList<Integer>"

# --- 3. Update the v1.1 sheet's data ---------------------------------------
# New (smaller/faster) raw timing samples for the ArrayList<Integer> version.
$ws2.Range("A2").Formula = "=3+122+39"
$ws2.Range("A3").Formula = "=2+406+142"
$ws2.Range("A4").Formula = "=7+1514+631"
$ws2.Range("A5").Formula = "=6+6322+3330"
$ws2.Range("A6").Formula = "=13+25761+11917"
$ws2.Range("A7").Formula = "=25+102032+48804"

# Only 6 size/time samples remain, so the two largest rows are now unused.
$ws2.Range("A8:B9").Clear()

# The FORECAST.LINEAR ranges now span just the 6 remaining samples (rows 2-7)
# instead of the original 8 (rows 2-9).
$ws2.Range("E2").Formula = "=_xlfn.FORECAST.LINEAR(D2,B`$2:B`$7,D`$2:D`$7)"
$ws2.Range("E3:E18").Formula = "=_xlfn.FORECAST.LINEAR(D3,B`$2:B`$7,D`$2:D`$7)"

# --- 4. Restore per-tab selections ------------------------------------------
$ws1.Select()
$ws1.Range("C34").Select()

$ws2.Select()
$ws2.Range("J24").Select()
